# Mostrar gráficos mensuales, semestrales y anuales
# Update "Estado" (status) column: all orders are now "Retirado" (picked up),
# replacing the former "Pendiente" (pending) / "Finalizado" (finished) values.
# Also push the estimated pickup date of the last order into November so the
# sample data spans more months for the new monthly/semestral/annual charts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "Estado" for rows 2-9.
$ws.Range("G2:G9").Value = "Retirado"

# Row 9's "Fecha de retiro estimado" (column C) moves from October to November.
$ws.Range("C9").Value = "20-11-2019"
